$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.167.86'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '2.053.88'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.54%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +3.70%  '
$ws.Range("E9").Value = '  +3.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").Value = '2.353.07'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.778'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '2.052.52'
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").Value = '37.078.11'
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.84%  '
$ws.Range("D22").Value = '0.0₃0807'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '224.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.50%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.18%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0615'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.70'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.475.90'
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("E45").Value = '  +4.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0930'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("E47").Value = '  +3.21%  '
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.09%  '
$ws.Range("E51").Value = '  +1.69%  '
